$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("MicroInstructions")

# Fix the dash typo in the COM instruction's semantic ("en dash" -> regular hyphen)
$ws.Range("B22").Value = "0xFF - DST "

# Add a colon after the condition in each conditional-branch semantic so the
# tokenizer can split "condition" from "action" (if X = Y PC = ... -> if X = Y: PC = ...)
$ws.Range("B60").Value = "if ccc = 0: PC = PC + IMM8s"
$ws.Range("B61").Value = "if ccc = 1: PC = PC + IMM8s"
$ws.Range("B62").Value = "if C = 0: PC = PC + IMM8s"
$ws.Range("B63").Value = "if C = 1: PC = PC + IMM8s"
$ws.Range("B64").Value = "if C = 0: PC = PC + IMM8s"
$ws.Range("B65").Value = "if C = 1: PC = PC + IMM8s"
$ws.Range("B66").Value = "if Z = 0: PC = PC + IMM8s"
$ws.Range("B67").Value = "if Z = 1: PC = PC + IMM8s"
$ws.Range("B68").Value = "if N = 0: PC = PC + IMM8s"
$ws.Range("B69").Value = "if N = 1: PC = PC + IMM8s"
$ws.Range("B70").Value = "if V = 0: PC = PC + IMM8s"
$ws.Range("B71").Value = "if V = 1: PC = PC + IMM8s"
$ws.Range("B72").Value = "if S = 0: PC = PC + IMM8s"
$ws.Range("B73").Value = "if S = 1: PC = PC + IMM8s"

# Leave the UI focused where the author ended up: MicroInstructions, cell B23 selected
$ws.Activate()
$ws.Range("B23").Select()
